$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column G: "Rodada 6" header, matching the style of the existing
# Rodada headers in B1:F1 (bold, centered, bordered) by copying the F1 format.
$ws.Range("G1").Value = "Rodada 6"
$ws.Range("F1").Copy()
$ws.Range("G1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Rodada 6 scores for each team (rows 2-68)
$ws.Range("G2").Value = 103.259765625
$ws.Range("G3").Value = 78.85986328125
$ws.Range("G4").Value = 80.06005859375
$ws.Range("G5").Value = 75.259765625
$ws.Range("G6").Value = 78.919921875
$ws.Range("G7").Value = 88.06005859375
$ws.Range("G8").Value = 102.06005859375
$ws.Range("G9").Value = 88.759765625
$ws.Range("G10").Value = 84.0498046875
$ws.Range("G11").Value = 92.4599609375
$ws.Range("G12").Value = 86.9599609375
$ws.Range("G13").Value = 82.60986328125
$ws.Range("G14").Value = 70.93017578125
$ws.Range("G15").Value = 68.4501953125
$ws.Range("G16").Value = 98.7998046875
$ws.Range("G17").Value = 107.66015625
$ws.Range("G18").Value = 99.56005859375
$ws.Range("G19").Value = 49.260009765625
$ws.Range("G20").Value = 98.2099609375
$ws.Range("G21").Value = 92.06005859375
$ws.Range("G22").Value = 102.85986328125
$ws.Range("G23").Value = 58.300048828125
$ws.Range("G24").Value = 102.85986328125
$ws.Range("G25").Value = 98.16015625
$ws.Range("G26").Value = 98.16015625
$ws.Range("G27").Value = 98.16015625
$ws.Range("G28").Value = 106.06005859375
$ws.Range("G29").Value = 77.5
$ws.Range("G30").Value = 98.2998046875
$ws.Range("G31").Value = 104.2001953125
$ws.Range("G32").Value = 86.06005859375
$ws.Range("G33").Value = 108.56005859375
$ws.Range("G34").Value = 48.840087890625
$ws.Range("G35").Value = 98.06005859375
$ws.Range("G36").Value = 109.64990234375
$ws.Range("G37").Value = 110.31982421875
$ws.Range("G38").Value = 106.14990234375
$ws.Range("G39").Value = 85.16015625
$ws.Range("G40").Value = 95.31982421875
$ws.Range("G41").Value = 66.759765625
$ws.Range("G42").Value = 103.66015625
$ws.Range("G43").Value = 119.85009765625
$ws.Range("G44").Value = 46.10009765625
$ws.Range("G45").Value = 131.6103515625
$ws.Range("G46").Value = 86.06005859375
$ws.Range("G47").Value = 110.35986328125
$ws.Range("G48").Value = 107.7001953125
$ws.Range("G49").Value = 106.93017578125
$ws.Range("G50").Value = 73.259765625
$ws.Range("G51").Value = 82.64990234375
$ws.Range("G52").Value = 86.35986328125
$ws.Range("G53").Value = 92.0498046875
$ws.Range("G54").Value = 78.56005859375
$ws.Range("G55").Value = 59.969970703125
$ws.Range("G56").Value = 111.43017578125
$ws.Range("G57").Value = 77.47998046875
$ws.Range("G58").Value = 88.35986328125
$ws.Range("G59").Value = 69.60986328125
$ws.Range("G60").Value = 100.85986328125
$ws.Range("G61").Value = 93.35986328125
$ws.Range("G62").Value = 62.43994140625
$ws.Range("G63").Value = 30.25
$ws.Range("G64").Value = 61.8798828125
$ws.Range("G65").Value = 77.85986328125
$ws.Range("G66").Value = 61.159912109375
$ws.Range("G67").Value = 65.5400390625
$ws.Range("G68").Value = 96.75
